$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.272198915481567
$ws.Range("B1").Value = 3.237287759780884
$ws.Range("C1").Value = 5.930813789367676
$ws.Range("D1").Value = 1.770566463470459
$ws.Range("E1").Value = 1.039753437042236
